# Update MSME country indicator figures for Madagascar Summary sheet.
# Source values gain extra decimal-place precision; cells are stored as
# text (shared strings) in the original workbook, so we force text entry
# with a leading apostrophe (quote-prefix) to avoid Excel's automatic
# "looks like a number" coercion, then restore the default "Normal" cell
# style so no stray NumberFormat/quote-prefix style is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$NewValue
    )
    $range = $ws.Range($Address)
    $range.Value = "'" + $NewValue
    $range.Style = "Normal"
}

# Enterprises density (per 1000 people) - row 11
Set-TextValue "B11" "11.13"
Set-TextValue "D11" "11.53"

# Employment (% of total) - row 12
Set-TextValue "C12" "18.23"
Set-TextValue "D12" "77.23"

# Enterprises (% of total) - row 14
Set-TextValue "B14" "96.34"
Set-TextValue "C14" "3.47"
Set-TextValue "D14" "99.81"
